$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing "Status" (column E) values moved up to fully-complete (1) for
# several tasks that had been left blank or at 90%.
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("E7").Value = 1
$ws.Range("E8").Value = 1

# Add a new "Code consolidation" task row, based on the last existing row
# (row 18) so it keeps the same shading/border/number formatting, then
# overwrite the two cells that actually differ (Output + Accountable).
$ws.Range("A18:F18").Copy($ws.Range("A19:F19"))
$ws.Range("A19").Value = "Code consolidation"
$ws.Range("D19").Value = "Kevin"

# Restore the active selection to what it was left at in the edited file.
$ws.Range("E6").Select()
